$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target table (players, positions, teams) for rows 2..19
$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Kel'el Ware", "C", "Miami Heat"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Jordan Hawkins", "SG", "New Orleans Pelicans"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
